# Update the "想去人数" (F column) counts on the 展览 and 全部类型 sheets.
# Both sheets carry the same event listing, so the same row/value updates
# apply to each.

$wb = $excel.ActiveWorkbook

# row -> new F-column value
$updates = @{
    2  = 17
    3  = 86
    4  = 259
    6  = 536
    7  = 50
    8  = 1988
    10 = 94
    11 = 4227
    12 = 33
    13 = 278
    15 = 96
    16 = 21
    17 = 58
    18 = 2963
    19 = 56
    20 = 416
    21 = 17
    22 = 14
    23 = 64
    28 = 45
    29 = 191
    30 = 6
    31 = 383
    32 = 1662
    33 = 238
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
